$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Disease_Synonymous")
$ws.Activate()

$ws.Range("A5").Value = "BO_H74"
$ws.Range("B5").Value = "H73"

$ws.Range("B6").Select()
